$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell $ws.Range("D2") "44.515.86"
Set-TextCell $ws.Range("E2") "  -5.03%  "

# Row 3
Set-TextCell $ws.Range("D3") "2.660.73"
Set-TextCell $ws.Range("E3") "  +1.52%  "

# Row 4
Set-TextCell $ws.Range("D4") "0.998"
Set-TextCell $ws.Range("E4") "  -0.02%  "

# Row 5
Set-TextCell $ws.Range("D5") "306.53"
Set-TextCell $ws.Range("E5") "  -0.50%  "

# Row 6
Set-TextCell $ws.Range("D6") "97.08"

# Row 7
Set-TextCell $ws.Range("D7") "0.589"
Set-TextCell $ws.Range("E7") "  -2.22%  "

# Row 8
Set-TextCell $ws.Range("E8") "  +0.01%  "

# Row 9
Set-TextCell $ws.Range("D9") "0.570"
Set-TextCell $ws.Range("E9") "  -3.05%  "

# Row 10
Set-TextCell $ws.Range("D10") "37.82"
Set-TextCell $ws.Range("E10") "  -4.88%  "

# Row 11
Set-TextCell $ws.Range("D11") "0.0831"
Set-TextCell $ws.Range("E11") "  -2.66%  "

# Row 12
Set-TextCell $ws.Range("D12") "7.97"
Set-TextCell $ws.Range("E12") "  -3.57%  "

# Row 13
Set-TextCell $ws.Range("D13") "3.052.02"
Set-TextCell $ws.Range("E13") "  +1.31%  "

# Row 14
Set-TextCell $ws.Range("E14") "  +0.85%  "

# Row 15
Set-TextCell $ws.Range("D15") "2.654.72"
Set-TextCell $ws.Range("E15") "  -2.47%  "

# Row 16
Set-TextCell $ws.Range("D16") "0.913"
Set-TextCell $ws.Range("E16") "  -1.32%  "

# Row 17
Set-TextCell $ws.Range("D17") "14.88"
Set-TextCell $ws.Range("E17") "  -1.60%  "

# Row 18
Set-TextCell $ws.Range("D18") "44.471.41"
Set-TextCell $ws.Range("E18") "  -5.28%  "

# Row 19
Set-TextCell $ws.Range("D19") "6.81"
Set-TextCell $ws.Range("E19") "  +1.51%  "

# Row 20
Set-TextCell $ws.Range("D20") "0.0₃0998"
Set-TextCell $ws.Range("E20") "  -2.15%  "

# Row 21
Set-TextCell $ws.Range("D21") "12.55"
Set-TextCell $ws.Range("E21") "  -4.21%  "

# Row 22
Set-TextCell $ws.Range("D22") "74.29"
Set-TextCell $ws.Range("E22") "  +2.53%  "

# Row 23
Set-TextCell $ws.Range("D23") "276.27"
Set-TextCell $ws.Range("E23") "  -0.22%  "

# Row 24
Set-TextCell $ws.Range("D24") "2.29"
Set-TextCell $ws.Range("E24") "  +5.08%  "

# Row 25
Set-TextCell $ws.Range("D25") "3.00"
Set-TextCell $ws.Range("E25") "  -1.42%  "

# Row 26
Set-TextCell $ws.Range("D26") "30.75"
Set-TextCell $ws.Range("E26") "  +0.67%  "

# Row 27
Set-TextCell $ws.Range("D27") "0.999"
Set-TextCell $ws.Range("E27") "  +0.04%  "

# Row 28
Set-TextCell $ws.Range("D28") "10.42"
Set-TextCell $ws.Range("E28") "  -2.07%  "

# Row 29
Set-TextCell $ws.Range("D29") "2.25"
Set-TextCell $ws.Range("E29") "  -2.82%  "

# Row 30
Set-TextCell $ws.Range("D30") "37.78"
Set-TextCell $ws.Range("E30") "  -4.38%  "

# Row 31
Set-TextCell $ws.Range("D31") "6.12"
Set-TextCell $ws.Range("E31") "  -1.19%  "

# Row 32
Set-TextCell $ws.Range("D32") "3.74"
Set-TextCell $ws.Range("E32") "  +3.03%  "

# Row 33
Set-TextCell $ws.Range("D33") "2.32"
Set-TextCell $ws.Range("E33") "  +5.01%  "

# Row 34
Set-TextCell $ws.Range("D34") "153.43"
Set-TextCell $ws.Range("E34") "  +2.11%  "

# Row 35
Set-TextCell $ws.Range("E35") "  -2.01%  "

# Row 36
Set-TextCell $ws.Range("D36") "0.0828"
Set-TextCell $ws.Range("E36") "  -2.91%  "

# Row 37
Set-TextCell $ws.Range("D37") "0.119"
Set-TextCell $ws.Range("E37") "  -3.48%  "

# Row 38
Set-TextCell $ws.Range("D38") "24.98"
Set-TextCell $ws.Range("E38") "  +8.99%  "

# Row 39
Set-TextCell $ws.Range("D39") "0.123"
Set-TextCell $ws.Range("E39") "  -0.58%  "

# Row 40
Set-TextCell $ws.Range("D40") "15.78"
Set-TextCell $ws.Range("E40") "  -1.52%  "

# Row 41
Set-TextCell $ws.Range("D41") "3.57"
Set-TextCell $ws.Range("E41") "  -2.11%  "

# Row 42
Set-TextCell $ws.Range("D42") "0.0321"
Set-TextCell $ws.Range("E42") "  -3.33%  "

# Row 43
Set-TextCell $ws.Range("B43") "RenderToken"
Set-TextCell $ws.Range("C43") "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell $ws.Range("D43") "3.92"
Set-TextCell $ws.Range("E43") "  -5.85%  "

# Row 44
Set-TextCell $ws.Range("B44") "Maker"
Set-TextCell $ws.Range("C44") "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextCell $ws.Range("D44") "2.114.91"
Set-TextCell $ws.Range("E44") "  -2.91%  "

# Row 45
Set-TextCell $ws.Range("D45") "0.996"
Set-TextCell $ws.Range("E45") "  -0.05%  "

# Row 46
Set-TextCell $ws.Range("D46") "91.76"
Set-TextCell $ws.Range("E46") "  -4.69%  "

# Row 47
Set-TextCell $ws.Range("D47") "9.36"
Set-TextCell $ws.Range("E47") "  -4.31%  "

# Row 48
Set-TextCell $ws.Range("B48") "Aave"
Set-TextCell $ws.Range("C48") "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextCell $ws.Range("D48") "109.60"
Set-TextCell $ws.Range("E48") "  -0.55%  "

# Row 49
Set-TextCell $ws.Range("B49") "RocketPoolETH"
Set-TextCell $ws.Range("C49") "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextCell $ws.Range("D49") "2.904.70"
Set-TextCell $ws.Range("E49") "  +1.38%  "

# Row 50
Set-TextCell $ws.Range("E50") "  +3.43%  "

# Row 51
Set-TextCell $ws.Range("D51") "0.196"
Set-TextCell $ws.Range("E51") "  -2.72%  "
